$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.743.89"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.627.60"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").Value = "'214.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").Value = "'0.0791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").Value = "'4.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "1.853.00"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.552"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.545.28"
$ws.Range("E15").Value = "  -6.71%  "
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").Value = "'62.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "25.729.61"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "'191.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").Value = "'6.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").Value = "'1.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("D26").Value = "'142.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("E27").Value = "  +3.36%  "
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").Value = "'15.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("D35").Value = "'2.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "'0.906"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").Value = "1.137.66"
$ws.Range("E38").Value = "  -2.26%  "
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'100.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "1.762.38"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "'55.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.41%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.417"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("E51").Value = "  -0.55%  "
